$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.060.08"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.593.66"
$ws.Range("E3").Value = "  +8.53%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'305.84"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'99.67"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +5.16%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.570"
$ws.Range("E9").Value = "  +11.55%  "
$ws.Range("D10").Value = "'38.43"
$ws.Range("E10").Value = "  +10.81%  "
$ws.Range("D11").Value = "'0.0832"
$ws.Range("E11").Value = "  +5.05%  "
$ws.Range("D12").Value = "'8.07"
$ws.Range("E12").Value = "  +12.83%  "
$ws.Range("D13").Value = "2.984.73"
$ws.Range("E13").Value = "  +8.56%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "2.605.63"
$ws.Range("E15").Value = "  +10.27%  "
$ws.Range("D16").Value = "'0.893"
$ws.Range("E16").Value = "  +8.24%  "
$ws.Range("D17").Value = "'14.81"
$ws.Range("E17").Value = "  +7.57%  "
$ws.Range("D18").Value = "46.203.73"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "'13.16"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("D20").Value = "0.0₃0999"
$ws.Range("E20").Value = "  +4.38%  "
$ws.Range("D21").Value = "'6.64"
$ws.Range("E21").Value = "  +9.04%  "
$ws.Range("D22").Value = "'70.65"
$ws.Range("E22").Value = "  +5.33%  "
$ws.Range("D23").Value = "'254.07"
$ws.Range("E23").Value = "  +3.84%  "
$ws.Range("D24").Value = "'2.97"
$ws.Range("E24").Value = "  +5.73%  "
$ws.Range("E25").Value = "  +14.88%  "
$ws.Range("D26").Value = "'27.86"
$ws.Range("E26").Value = "  +31.95%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'10.40"
$ws.Range("E28").Value = "  +6.09%  "
$ws.Range("D29").Value = "'39.69"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "'2.26"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("E31").Value = "  +9.74%  "
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'2.30"
$ws.Range("E33").Value = "  +17.95%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.92"
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("D35").Value = "'151.74"
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("D36").Value = "'0.0830"
$ws.Range("E36").Value = "  +6.89%  "
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("D39").Value = "'4.17"
$ws.Range("E39").Value = "  +6.13%  "
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("D41").Value = "'3.58"
$ws.Range("E41").Value = "  +10.22%  "
$ws.Range("E42").Value = "  +6.60%  "
$ws.Range("D43").Value = "2.042.07"
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("D44").Value = "'18.91"
$ws.Range("E44").Value = "  +32.80%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'90.77"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "'9.12"
$ws.Range("E47").Value = "  +7.51%  "
$ws.Range("D48").Value = "'109.33"
$ws.Range("E48").Value = "  +10.09%  "
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").Value = "2.843.89"
$ws.Range("E50").Value = "  +8.55%  "
$ws.Range("D51").Value = "'0.199"
$ws.Range("E51").Value = "  +6.54%  "
